# LinearizerAbsTaylor: switch column G ("Simulated-annealing-joa-parallel")
# from the SA.v7 run to the SA.v1 run on both the "Tiempos" and "Nodos"
# sheets, and move the conditional "no errors" highlight rules that used
# to track B5/C6/E8 so that they now track the newly-populated G5/G6/G8
# cells (rows 11 and 12 previously had no recorded value for this column
# at all).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tiempos")
$ws2 = $wb.Worksheets.Item("Nodos")

# --- Column G values: "Tiempos" sheet (SA.v7 -> SA.v1 timings) ---------
$tiemposG = @{
    "G2"  = 3.23
    "G3"  = 75.26000000000001
    "G4"  = 16.23
    "G5"  = 12.5
    "G6"  = 13.03
    "G7"  = 23.54
    "G8"  = 1.11
    "G9"  = 153.04
    "G10" = 54.18
    "G11" = 3.99
    "G12" = 2.93
    "G13" = 1.16
    "G14" = 87.01000000000001
}

foreach ($addr in $tiemposG.Keys) {
    $ws1.Range($addr).Value2 = $tiemposG[$addr]
}

# --- Column G values: "Nodos" sheet (SA.v7 -> SA.v1 node counts) -------
$nodosG = @{
    "G2"  = 632.8
    "G3"  = 12119.6
    "G4"  = 1653.6
    "G5"  = 5108.4
    "G6"  = 5081.6
    "G7"  = 11185.2
    "G8"  = 239.6
    "G9"  = 60294.8
    "G10" = 31781.2
    "G11" = 500
    "G12" = 999.2
    "G13" = 153.6
    "G14" = 2771.6
}

foreach ($addr in $nodosG.Keys) {
    $ws2.Range($addr).Value2 = $nodosG[$addr]
}

# --- Conditional formatting on "Nodos": relocate the notContainsErrors
# rules that applied to B5, C6, E8 onto the now-meaningful G5, G6, G8
# cells (same dxf/priority, new formula + AppliesTo range).
$moves = @(
    @{ From = "B5"; To = "G5" },
    @{ From = "C6"; To = "G6" },
    @{ From = "E8"; To = "G8" }
)

foreach ($mv in $moves) {
    $rule = $ws2.Range($mv.From).FormatConditions.Item(1)
    $rule.Formula1 = "NOT(ISERROR($($mv.To)))"
    $rule.ModifyAppliesToRange($ws2.Range($mv.To))
}
